$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing hours/work-type entries for Week 2 (rows 9 and 11)
$ws.Range("F9").Value = 1.5
$ws.Range("G9").Value = "Implementation"

$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = "Management"

# Nudge the window tab-split ratio to match the saved view state
$excel.ActiveWindow.TabRatio = 0.985

# Move the active selection to F9
$ws.Range("F9").Select() | Out-Null
